$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Atividades imobiliárias"
$ws.Range("B2").Value = 28.31208085602502
$ws.Range("C2").Value = "2014 / 2023"

# Row 3
$ws.Range("A3").Value = "Atividades financeiras, de seguros e serviços relacionados"
$ws.Range("B3").Value = 23.65397356444598
$ws.Range("C3").Value = "2014 / 2023"

# Row 4
$ws.Range("A4").Value = "Eletricidade e gás, água, esgoto, atividades de gestão de resíduos e descontaminação"
$ws.Range("B4").Value = 22.70215620370851
$ws.Range("C4").Value = "2014 / 2023"

# Row 5
$ws.Range("A5").Value = "Informação e comunicação"
$ws.Range("B5").Value = 18.69747267971518
$ws.Range("C5").Value = "2014 / 2023"

# Row 6
$ws.Range("A6").Value = "Agropecuária"
$ws.Range("B6").Value = 5.682545327113388
$ws.Range("C6").Value = "2014 / 2023"

# Row 7
$ws.Range("A7").Value = "Administração, defesa, educação e saúde públicas e seguridade social"
$ws.Range("B7").Value = 2.548568009332001
$ws.Range("C7").Value = "2014 / 2023"

# Row 8
$ws.Range("A8").Value = "Indústrias extrativas"
$ws.Range("B8").Value = 69.71390886407463
$ws.Range("C8").Value = "2022 / 2023"

# Row 9
$ws.Range("A9").Value = "Agropecuária"
$ws.Range("B9").Value = 7.64923808343667
$ws.Range("C9").Value = "2022 / 2023"

# Row 10
$ws.Range("A10").Value = "Informação e comunicação"
$ws.Range("B10").Value = 7.133426276710783
$ws.Range("C10").Value = "2022 / 2023"

# Row 11
$ws.Range("A11").Value = "Atividades financeiras, de seguros e serviços relacionados"
$ws.Range("B11").Value = 4.80083445372183
$ws.Range("C11").Value = "2022 / 2023"

# Row 12
$ws.Range("A12").Value = "Comércio e reparação de veículos automotores e motocicletas"
$ws.Range("B12").Value = 2.010492569309139
$ws.Range("C12").Value = "2022 / 2023"

# Row 13
$ws.Range("A13").Value = "Indústrias de transformação"
$ws.Range("B13").Value = 1.925088415038929
$ws.Range("C13").Value = "2022 / 2023"
